$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata fixes ---
$ws.Range("I2").Value = "Gropin growth model for Aeromonas hydrophila in/on modified BHI (gropin ID: 24 )"
$ws.Range("I7").Value = "15/10/2020"

# --- Rows 133-135: description + uppercase DOUBLE, drop source/subject/dist ---
$descEnvFactor = "Only applicaple in either mode 'responsesurface' or 'time2multiply'. Environmental factors for growth model."
foreach ($r in 133..135) {
    $ws.Range("O$r").Value = $descEnvFactor
    $ws.Range("R$r").Value = "DOUBLE"
    $ws.Range("S$r").ClearContents()
    $ws.Range("T$r").ClearContents()
    $ws.Range("U$r").ClearContents()
}

# --- Rows 136-139: drop description/source/subject/dist (R stays "double") ---
foreach ($r in 136..139) {
    $ws.Range("O$r").ClearContents()
    $ws.Range("S$r").ClearContents()
    $ws.Range("T$r").ClearContents()
    $ws.Range("U$r").ClearContents()
}

# --- Row 140-141: visualisation axis description, drop source/subject/dist ---
$descVisAxis = "For visualisation purposes in either mode 'time2multiply' or 'responsesurface'. visualisation axis. Enter string with '<variable ID>'. Strings that are accepted: T, aw, CO2dissolved"
foreach ($r in 140..141) {
    $ws.Range("O$r").Value = $descVisAxis
    $ws.Range("S$r").ClearContents()
    $ws.Range("T$r").ClearContents()
    $ws.Range("U$r").ClearContents()
}

# --- Row 142: mode description ---
$ws.Range("O142").Value = "three different modes are available: 'responsesurface' is running the secondary model calculating mumax only. 'time2multiply' returns a 2D-plot of the time the microorganism needs to increase N by a logstep of 'logIncrease'(free parameter to choose). 'kinetic' runs the tertiary model, based on the variables chosen (with '_kinetic'-suffix)."
$ws.Range("S142").ClearContents()
$ws.Range("T142").ClearContents()
$ws.Range("U142").ClearContents()

# --- Row 143: lagTime description ---
$ws.Range("O143").Value = "Only applicaple in either mode 'time2multiply' or 'kinetic' is chosen. This is the time the microorganism needs for adjusting to its environment before multiplying."
$ws.Range("S143").ClearContents()
$ws.Range("T143").ClearContents()
$ws.Range("U143").ClearContents()

# --- Row 144: logIncrease description ---
$ws.Range("O144").Value = "Only applicaple in mode 'time2multiply'. Free parameter to calculate the time the microorganism needs to increase its numbers by the log step increase indictated by this value."
$ws.Range("S144").ClearContents()
$ws.Range("T144").ClearContents()
$ws.Range("U144").ClearContents()

# --- Row 145: logN0 description ---
$ws.Range("O145").Value = "Only applicaple in either mode 'kinetic'. Choose the number of microorganisms at the beginning of this simulation. (log step!)"
$ws.Range("S145").ClearContents()
$ws.Range("T145").ClearContents()
$ws.Range("U145").ClearContents()

# --- Row 146: logNEnd description ---
$ws.Range("O146").Value = "Only applicaple in either mode 'kinetic'. Choose the number of microorganisms at the end of this simulation. (log step!)"
$ws.Range("S146").ClearContents()
$ws.Range("T146").ClearContents()
$ws.Range("U146").ClearContents()

# --- Row 147: simTime description ---
$ws.Range("O147").Value = "time of simulation, unit is h-1"
$ws.Range("S147").ClearContents()
$ws.Range("T147").ClearContents()
$ws.Range("U147").ClearContents()

# --- Rows 148-150: kinetic variable description (shared text) ---
$descKineticVar = "Only applicaple in either mode 'kinetic'. Choose variable for prediction of growth depending on environmental factors."
foreach ($r in 148..150) {
    $ws.Range("O$r").Value = $descKineticVar
    $ws.Range("S$r").ClearContents()
    $ws.Range("T$r").ClearContents()
    $ws.Range("U$r").ClearContents()
}
